# Add "Premium Cat Food" to Beth S.'s client list (sheet1), matching the
# existing shared string used already on Chris K.'s sheet, then make
# Beth S. the active/selected sheet (instead of Mary M.).

$wb = $excel.ActiveWorkbook

$wsBeth = $wb.Worksheets.Item("Beth S.")

# Append the new item in the first empty row of column A.
$wsBeth.Range("A6").Value = "Premium Cat Food"

# Select the newly-added cell on Beth S.'s sheet.
$wsBeth.Range("A6").Select()

# Make Beth S. the active sheet (was Mary M.) and deselect Mary M. as the tab.
$wsBeth.Activate()

$wb.Save()
